$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42607.89025462963
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 55
$ws.Range("D5").Value = 42
$ws.Range("E5").Value = 54
$ws.Range("F5").Value = 45
$ws.Range("G5").Value = 13120
$ws.Range("H5").Value = 23309
$ws.Range("I5").Value = 2718
$ws.Range("J5").Value = 292
$ws.Range("K5").Value = 226
$ws.Range("L5").Value = 12
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = "Bag"
